# Applies the "Updated symbol list" data refresh (Mon Jan 16 15:15:57 UTC 2023).
# Every value in this sheet is stored as literal text (t="inlineStr" in the source
# OOXML), including numeric-looking price/volume/hour columns, so values are written
# with a leading apostrophe to force text entry, then ClearFormats() removes the
# auto-applied "@" text NumberFormat/style so the cell keeps the workbook's original
# (unstyled) look - matching how the source file had no per-cell style on data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ "D2"="294.73"; "E2"="-2.20%"; "G2"="15" }
    @{ "D3"="31.31"; "E3"="-0.59%"; "G3"="15" }
    @{ "D4"="5.088"; "E4"="-1.24%"; "G4"="15" }
    @{ "E5"="8.59%"; "G5"="15" }
    @{ "D6"="2.489"; "E6"="34.98%"; "G6"="15" }
    @{ "D7"="7.760"; "E7"="-0.43%"; "G7"="15" }
    @{ "B8"="GateToken"; "C8"="https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"; "D8"="3.804"; "E8"="1.42%"; "G8"="15" }
    @{ "B9"="MXToken"; "C9"="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; "D9"="0.9252"; "E9"="-0.18%"; "G9"="15" }
    @{ "B10"="WazirX"; "C10"="https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; "D10"="0.1788"; "E10"="6.05%"; "G10"="15" }
    @{ "B11"="LiechtensteinCryptoassetsExchange"; "C11"="https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"; "D11"="0.07352"; "E11"="4.72%"; "G11"="15" }
    @{ "B12"="MandalaExchangeToken"; "C12"="https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; "D12"="0.09144"; "E12"="13.61%"; "G12"="15" }
    @{ "B13"="BitrueCoin"; "C13"="https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; "D13"="0.03037"; "E13"="0.07%"; "G13"="15" }
    @{ "B14"="BitMartToken"; "C14"="https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; "D14"="0.09990"; "E14"="0.55%"; "G14"="15" }
    @{ "B15"="BitForexToken"; "C15"="https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; "D15"="0.001496"; "E15"="-0.21%"; "G15"="15" }
    @{ "B16"="TigerCash"; "C16"="https://coinranking.com/coin/6hIn06L2+tigercash-tch"; "D16"="0.005925"; "E16"="-3.38%"; "G16"="15" }
    @{ "B17"="LEO"; "C17"="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; "D17"="3.509"; "E17"="1.39%"; "G17"="15" }
    @{ "D18"="2.247"; "E18"="1.03%"; "G18"="15" }
    @{ "E19"="-0.65%"; "G19"="15" }
    @{ "D20"="0.1338"; "G20"="15" }
    @{ "D21"="4.300"; "E21"="-5.47%"; "G21"="15" }
    @{ "D22"="0.1619"; "E22"="2.34%"; "G22"="15" }
    @{ "D23"="0.04599"; "E23"="-0.70%"; "G23"="15" }
    @{ "D24"="0.001251"; "E24"="3.01%"; "G24"="15" }
    @{ "D25"="0.004399"; "E25"="-7.44%"; "G25"="15" }
    @{ "E26"="-7.29%"; "G26"="15" }
    @{ "D27"="0.0003432"; "E27"="83.37%"; "G27"="15" }
    @{ "G28"="15" }
    @{ "G29"="15" }
    @{ "G30"="15" }
    @{ "G31"="15" }
    @{ "G32"="15" }
    @{ "G33"="15" }
    @{ "G34"="15" }
    @{ "G35"="15" }
    @{ "G36"="15" }
    @{ "G37"="15" }
    @{ "G38"="15" }
    @{ "D39"="0.01746"; "E39"="1.50%"; "G39"="15" }
    @{ "D40"="0.04434"; "E40"="-1.30%"; "G40"="15" }
    @{ "D41"="0.006880"; "E41"="-2.84%"; "G41"="15" }
    @{ "D42"="0.1338"; "E42"="-0.30%"; "G42"="15" }
    @{ "D43"="0.002148"; "E43"="-0.71%"; "G43"="15" }
    @{ "D44"="0.009769"; "E44"="-10.86%"; "G44"="15" }
    @{ "D45"="0.00006581"; "E45"="5.48%"; "G45"="15" }
    @{ "E46"="0.13%"; "G46"="15" }
    @{ "G47"="15" }
    @{ "E48"="-55.53%"; "G48"="15" }
    @{ "D49"="0.00002103"; "E49"="0.13%"; "G49"="15" }
    @{ "D50"="0.0002003"; "E50"="0.20%"; "G50"="15" }
    @{ "G51"="15" }
)

foreach ($rowUpdate in $updates) {
    foreach ($cellRef in $rowUpdate.Keys) {
        $range = $ws.Range($cellRef)
        $range.Value = "'" + $rowUpdate[$cellRef]
        $range.ClearFormats()
    }
}
